# "adding div class to excel"
#
# The sheet lists numbered rows (A: 1-5) with a short word in column B
# (Amplify, Revolutionize, Empower, Optimize, Elevate). This adds a third
# column, C, holding the corresponding remixicon CSS class name to use for
# each row's icon ("div class=...") on the welcome page.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for column C, one per existing row (rows 1-5).
$iconClasses = @(
    "ri-store-line",
    "ri-bar-chart-box-line",
    "ri-calendar-todo-line",
    "ri-paint-brush-line",
    "ri-database-2-line"
)

for ($i = 0; $i -lt $iconClasses.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 3).Value = $iconClasses[$i]
}

# Let the sheet re-fit column widths (B got narrower, C is sized to its
# new content) now that a column has been added.
$ws.Columns.Item(2).EntireColumn.AutoFit() | Out-Null
$ws.Columns.Item(3).EntireColumn.AutoFit() | Out-Null

# Reset the selection back to the top-left cell instead of the whole
# first row that was selected before.
$ws.Range("A1").Select()
